# "Rajout de la partie d'Anis sur la DCT"
#
# - Salaires - Salaires : Thomas Debarre salaire brut passe de 2750 a 2700
# - Compte de resultat previsionnel :
#     * I3 (Licences, deja deduit) devient une formule vers le total des
#       salaires au lieu d'une constante
#     * F4 (nombre d'unites du Logiciel) passe de 24000 a 26000
#     * Nouvelle ligne de charges deja deduites "2 serveurs" (H5/I5)
# Les autres feuilles (Impots, Bilan - Tableau 1, Bilan - Tableau 2) ne
# contiennent que des formules qui se recalculent automatiquement.

$wb = $excel.ActiveWorkbook

# --- Salaires - Salaires ---------------------------------------------------
$wsSalaires = $wb.Worksheets.Item("Salaires - Salaires")
$wsSalaires.Range("C17").Value = 2700

# --- Compte de resultat previsionnel ---------------------------------------
$wsCompte = $wb.Worksheets.Item("Compte de résultat prévisionnel")

# Le total des charges salariales (precedemment une valeur figee) pointe
# maintenant vers le total de la feuille Salaires.
$wsCompte.Range("I3").Formula = "='Salaires - Salaires'!G18"

# Nombre de licences logiciel vendues : 24000 -> 26000
$wsCompte.Range("F4").Value = 26000

# Nouvelle ligne "2 serveurs" (partie d'Anis sur la DCT)
# On reprend la mise en forme de H4 (meme "libelle" sur la ligne du dessus)
$wsCompte.Range("H4").Copy()
$wsCompte.Range("H5").PasteSpecial(-4122) # xlPasteFormats
$wsCompte.Range("H5").Value = "2 serveurs"
$wsCompte.Range("I5").Formula = "=12*606"

$wb.Application.CalculateFullRebuild()
